$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.341.72'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '3.332.85'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.90'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.75%  '
$ws.Range('E7').Value = '  +1.53%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '3.323.12'
$ws.Range('E9').Value = '  -0.99%  '
$ws.Range('E10').Value = '  +6.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.636'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.01%  '
$ws.Range('E13').Value = '  +1.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.06'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').Value = '3.863.84'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('E16').Value = '  +2.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.10'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').Value = '3.334.24'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '64.137.28'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.985'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '447.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.97'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('E24').Value = '  -1.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.54'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.73'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.62'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.59%  '
$ws.Range('E29').Value = '  -2.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.82'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.50'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '62.21'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '568.17'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('E35').Value = '  -1.30%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.51'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.25'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('D41').Value = '0.0₃0729'
$ws.Range('E41').Value = '  -3.79%  '
$ws.Range('D42').Value = '3.058.27'
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('E44').Value = '  -3.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.134'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.03%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.38%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.07%  '
$ws.Range('E51').Value = '  -1.63%  '
